$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 edits: "fda" -> "fdafd", "fda" -> "fd", "2" -> "1"
$ws.Range("A2").Value = "fdafd"
$ws.Range("C2").Value = "fd"

# D2 needs to stay a text value ("1"), not get auto-converted to a number.
# Temporarily format as text so COM keeps it a string, then restore the
# cell's style back to Normal so no stray number-format style lingers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "1"
$ws.Range("D2").Style = "Normal"

# Row 3 was merged into A2 (A2:A3); unmerge before removing the row.
$ws.Range("A2:A3").UnMerge()

# Remove row 3 entirely (its data is gone in the new version).
$ws.Rows(3).Delete()
